$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamps for the existing rows (C2:C4)
# Old value 45759.91172538449 -> new value 45759.91172538194
$ws.Range("C2:C4").Value = 45759.91172538194

# Append new rows for an updated quotation snapshot
$ws.Range("A5").Value = "Dólar"
$ws.Range("B5").Value = 5.8546
$ws.Range("C5").Value = 45761.82636843352

$ws.Range("A6").Value = "Euro"
$ws.Range("B6").Value = 6.64011
$ws.Range("C6").Value = 45761.82636843352

$ws.Range("A7").Value = "Bitcoin"
$ws.Range("B7").Value = 497501000
$ws.Range("C7").Value = 45761.82636843352

# Apply the same number format used by C2:C4 (datetime format) to the new cells
$ws.Range("C5:C7").NumberFormat = $ws.Range("C4").NumberFormat
